# Updated cryptos list on Mon Apr 24 03:36:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: assign directly so they stay as text,
# matching the original inline-string cell type.
$textUpdates = @(
    @{Cell='D2'; Value='28.076.29'}
    @{Cell='E2'; Value='  +1.47%  '}
    @{Cell='D3'; Value='1.889.17'}
    @{Cell='E3'; Value='  +0.96%  '}
    @{Cell='E4'; Value='  +1.19%  '}
    @{Cell='E5'; Value='  +1.51%  '}
    @{Cell='E6'; Value='  +1.08%  '}
    @{Cell='E7'; Value='  +1.20%  '}
    @{Cell='E8'; Value='  +0.44%  '}
    @{Cell='E9'; Value='  -0.97%  '}
    @{Cell='E10'; Value='  -0.10%  '}
    @{Cell='E11'; Value='  -0.04%  '}
    @{Cell='E12'; Value='  +0.80%  '}
    @{Cell='D13'; Value='1.873.07'}
    @{Cell='E13'; Value='  +0.06%  '}
    @{Cell='E14'; Value='  +1.67%  '}
    @{Cell='E15'; Value='  +1.15%  '}
    @{Cell='E16'; Value='  +1.09%  '}
    @{Cell='E17'; Value='  +2.18%  '}
    @{Cell='E18'; Value='  +2.08%  '}
    @{Cell='E19'; Value='  +0.69%  '}
    @{Cell='E20'; Value='  -0.86%  '}
    @{Cell='E21'; Value='  +1.10%  '}
    @{Cell='D22'; Value='28.048.30'}
    @{Cell='E22'; Value='  +1.29%  '}
    @{Cell='E23'; Value='  +0.85%  '}
    @{Cell='E24'; Value='  +0.34%  '}
    @{Cell='E25'; Value='  +1.72%  '}
    @{Cell='D26'; Value='2.108.53'}
    @{Cell='E27'; Value='  +1.68%  '}
    @{Cell='E28'; Value='  -1.20%  '}
    @{Cell='E29'; Value='  +0.76%  '}
    @{Cell='E30'; Value='  -0.26%  '}
    @{Cell='E31'; Value='  -0.33%  '}
    @{Cell='E32'; Value='  +1.37%  '}
    @{Cell='E33'; Value='  +1.13%  '}
    @{Cell='E34'; Value='  +1.44%  '}
    @{Cell='E35'; Value='  +0.92%  '}
    @{Cell='E36'; Value='  -5.44%  '}
    @{Cell='B37'; Value='VeChain'}
    @{Cell='C37'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'}
    @{Cell='E37'; Value='  -0.16%  '}
    @{Cell='B38'; Value='Hedera'}
    @{Cell='C38'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'}
    @{Cell='E39'; Value='  -2.63%  '}
    @{Cell='E40'; Value='  +0.91%  '}
    @{Cell='E41'; Value='  +1.18%  '}
    @{Cell='E42'; Value='  -0.02%  '}
    @{Cell='E43'; Value='  +0.06%  '}
    @{Cell='E44'; Value='  +0.99%  '}
    @{Cell='E45'; Value='  +1.00%  '}
    @{Cell='E46'; Value='  -0.31%  '}
    @{Cell='E47'; Value='  -0.35%  '}
    @{Cell='E48'; Value='  -0.25%  '}
    @{Cell='E49'; Value='  -0.72%  '}
    @{Cell='E50'; Value='  -0.35%  '}
    @{Cell='E51'; Value='  -1.57%  '}
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Price values that look numeric (e.g. '1.015', '0.06082'): force the cell's
# number format to Text first so Excel stores the exact original string instead
# of silently converting it to a floating point number.
$numericLookingUpdates = @(
    @{Cell='D4'; Value='1.015'}
    @{Cell='D5'; Value='336.66'}
    @{Cell='D7'; Value='0.4747'}
    @{Cell='D8'; Value='0.3957'}
    @{Cell='D9'; Value='47.25'}
    @{Cell='D10'; Value='0.08042'}
    @{Cell='D11'; Value='1.022'}
    @{Cell='D12'; Value='21.96'}
    @{Cell='D14'; Value='6.038'}
    @{Cell='D15'; Value='7.226'}
    @{Cell='D16'; Value='1.016'}
    @{Cell='D17'; Value='88.57'}
    @{Cell='D18'; Value='0.06768'}
    @{Cell='D19'; Value='0.00001054'}
    @{Cell='D20'; Value='17.08'}
    @{Cell='D21'; Value='1.014'}
    @{Cell='D23'; Value='5.534'}
    @{Cell='D24'; Value='11.03'}
    @{Cell='D27'; Value='160.98'}
    @{Cell='D28'; Value='20.02'}
    @{Cell='D29'; Value='2.109'}
    @{Cell='D30'; Value='5.542'}
    @{Cell='D31'; Value='122.03'}
    @{Cell='D32'; Value='0.9792'}
    @{Cell='D33'; Value='0.09598'}
    @{Cell='D34'; Value='3.645'}
    @{Cell='D35'; Value='5.366'}
    @{Cell='D37'; Value='0.02255'}
    @{Cell='D38'; Value='0.06082'}
    @{Cell='D39'; Value='1.202'}
    @{Cell='D40'; Value='8.206'}
    @{Cell='D42'; Value='0.5976'}
    @{Cell='D43'; Value='0.1896'}
    @{Cell='D44'; Value='10.36'}
    @{Cell='D45'; Value='1.265'}
    @{Cell='D46'; Value='0.5670'}
    @{Cell='D47'; Value='12.14'}
    @{Cell='D48'; Value='1.933'}
    @{Cell='D49'; Value='3.366'}
    @{Cell='D50'; Value='0.06828'}
    @{Cell='D51'; Value='112.39'}
)

foreach ($u in $numericLookingUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
}

